$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ticker symbols for columns B (col 2) and C (col 3), rows 2-16
$colB = @(
    "NSE:ADROITINFO",
    "NSE:BDL",
    "NSE:FIEMIND",
    "NSE:HNGSNGBEES",
    "NSE:IMPAL",
    "NSE:JMFINANCIL",
    "NSE:JYOTHYLAB",
    "NSE:MAANALU",
    "NSE:NSLNISP",
    "NSE:PUNJABCHEM",
    "NSE:RAMRAT",
    $null,
    $null,
    $null,
    $null
)

$colC = @(
    "NSE:ADORWELD",
    "NSE:ASHIANA",
    "NSE:ASMS",
    "NSE:ASPINWALL",
    "NSE:BHARATFORG",
    "NSE:BLISSGVS",
    "NSE:EXXARO",
    "NSE:HCC",
    "NSE:HEIDELBERG",
    "NSE:HESTERBIO",
    "NSE:JBMA",
    "NSE:PATINTLOG",
    "NSE:PILANIINVS",
    "NSE:PRSMJOHNSN",
    "NSE:RAJRATAN"
)

for ($i = 0; $i -lt 15; $i++) {
    $row = 2 + $i

    if ($colB[$i]) {
        $ws.Cells.Item($row, 2).Value = $colB[$i]
    } else {
        $ws.Cells.Item($row, 2).Value = ""
    }

    $ws.Cells.Item($row, 3).Value = $colC[$i]
}

# Remove the now-unused rows 17-25 (sheet shrinks from 25 rows to 16 rows)
$ws.Range("A17:F25").EntireRow.Delete()
